# Update clp format note: clarify that -1 means no keyframe and should be skipped
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F9").Value = "repeat for every key frame / if -1 no keyframe and skip to next"

# Widen column F so the longer note text is readable
$ws.Columns("F").ColumnWidth = 57.16667

# Move the active selection to E14
$excel.Goto($ws.Range("E14"))
